# "this added by last report 15-05-25"
# Update the Route Cost RSO workbook's "Route" sheet:
#   - Report date (L3) changed from a date serial to the literal text "15/5/2025"
#   - Route cost (D7/D9/D10) figures revised, flowing through the I/L formula columns
#     and the Total row (12)
#   - Signature cell L17 ("Chiranjit Barai") cleared out

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Route")

# Report date, typed as plain text (day=15 doesn't parse as a valid date so it
# stays text, matching the shared-string cell in the target file)
$ws.Range("L3").Value = "15/5/2025"

# Revised route-cost figures
$ws.Range("D7").Value = 130
$ws.Range("D9").Value = 130
$ws.Range("D10").Value = 200

# Clear the stray name left in the signature row
$ws.Range("L17").ClearContents()
